# Append newly-scraped bike rows to Sheet2 (ZigWheels bike listing data).
# Row 2 used to be "Yamaha EC-06" / row 8 ended the table; the refreshed
# scrape shifts everything up by one (the EC-06 row drops off) and adds a
# further batch of freshly scraped bikes through row 29.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$rows = @(
    @(2,  "Yamaha AEROX-E",                 "Rs. 2.90 Lakh", "Expected Launch : Mar 2026"),
    @(3,  "2026 Yamaha MT-03",               "Rs. 3.60 Lakh", "Expected Launch : Jun 2026"),
    @(4,  "Yamaha NMax 155",                 "Rs. 1.40 Lakh", "Expected Launch : Jun 2026"),
    @(5,  "2026 Yamaha R3",                  "Rs. 3.70 Lakh", "Expected Launch : Jun 2026"),
    @(6,  "Yamaha WR155 R",                  "Rs. 1.70 Lakh", "Expected Launch : Jun 2026"),
    @(7,  "Yamaha RX 100",                   "Rs. 1.00 Lakh", "Expected Launch : Dec 2026"),
    @(8,  "BMW F 450 GS",                    "Rs. 4.50 Lakh", "Expected Launch : Mar 2026"),
    @(9,  "Hero Karizma XMR 250",            "Rs. 2.00 Lakh", "Expected Launch : Mar 2026"),
    @(10, "2026 Husqvarna Svartpilen 401",   "Rs. 2.99 Lakh", "Expected Launch : Mar 2026"),
    @(11, "2026 Husqvarna Vitpilen 250",     "Rs. 2.30 Lakh", "Expected Launch : Mar 2026"),
    @(12, "Royal Enfield Flying Flea C6",    "Rs. 2.00 Lakh", "Expected Launch : Mar 2026"),
    @(13, "Kawasaki W230",                   "Rs. 1.50 Lakh", "Expected Launch : Apr 2026"),
    @(14, "2026 Yamaha MT-03",               "Rs. 3.60 Lakh", "Expected Launch : Jun 2026"),
    @(15, "2026 Yezdi Scrambler",            "Rs. 2.15 Lakh", "Expected Launch : Jun 2026"),
    @(16, "BSA Scrambler",                   "Rs. 3.45 Lakh", "Expected Launch : Jun 2026"),
    @(17, "BSA Electric Bike",               "Rs. 2.50 Lakh", "Expected Launch : Jun 2026"),
    @(18, "Triumph Bonneville 350",          "Rs. 1.85 Lakh", "Expected Launch : Jun 2026"),
    @(19, "2026 KTM RC 390",                 "Rs. 3.50 Lakh", "Expected Launch : Jun 2026"),
    @(20, "2026 Yamaha R3",                  "Rs. 3.70 Lakh", "Expected Launch : Jun 2026"),
    @(21, "Yamaha WR155 R",                  "Rs. 1.70 Lakh", "Expected Launch : Jun 2026"),
    @(22, "Bajaj Platina 125",               "Rs. 80,000",    "Expected Launch : Jun 2026"),
    @(23, "Bajaj Pulsar NS150",              "Rs. 1.22 Lakh", "Expected Launch : Jun 2026"),
    @(24, "BSA Bantam 350",                  "Rs. 2.20 Lakh", "Expected Launch : Jul 2026"),
    @(25, "KTM 350 Duke",                    "Rs. 2.60 Lakh", "Expected Launch : Jul 2026"),
    @(26, "Ola Cruiser",                     "Rs. 2.70 Lakh", "Expected Launch : Aug 2026"),
    @(27, "Royal Enfield 250",               "Rs. 1.30 Lakh", "Expected Launch : Sep 2026"),
    @(28, "Royal Enfield Interceptor 750",   "Rs. 3.80 Lakh", "Expected Launch : Nov 2026"),
    @(29, "KTM Duke 490",                    "Rs. 4.00 Lakh", "Expected Launch : Nov 2026")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
